$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value  = 1.72
$ws.Range("I2").Value  = 5.7
$ws.Range("J2").Value  = 4.2
$ws.Range("K2").Value  = 4.3
$ws.Range("L2").Value  = 1.37
$ws.Range("N2").Value  = 4
$ws.Range("O2").Value  = 1.29
$ws.Range("V2").Value  = 1.21
$ws.Range("X2").Value  = 17.5
$ws.Range("Y2").Value  = 21
$ws.Range("Z2").Value  = 44
$ws.Range("AA2").Value = 150
$ws.Range("AC2").Value = 9.6
$ws.Range("AD2").Value = 21
$ws.Range("AE2").Value = 75
$ws.Range("AF2").Value = 10.5
$ws.Range("AG2").Value = 9.8
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 75
$ws.Range("AJ2").Value = 17
$ws.Range("AK2").Value = 17.5
$ws.Range("AL2").Value = 34
$ws.Range("AM2").Value = 110
$ws.Range("AO2").Value = 85

# Row 4 updates
$ws.Range("I4").Value  = 1.44
$ws.Range("J4").Value  = 4.6
$ws.Range("P4").Value  = 1.98
$ws.Range("U4").Value  = 1.74
$ws.Range("V4").Value  = 3.15
$ws.Range("AB4").Value = 30
$ws.Range("AO4").Value = 7.6
